# Saldo_guide.xlsx update
# - Rename sheet to reflect the new extraction run (2024-08-27 09:46:17)
# - Bump the "Dt. Referencia" column (G) from 45530 (2024-08-26) to 45531 (2024-08-27)
#   for every data row (2 through 274)
# - Correct a handful of projected/total balance values (columns E and H)
#   that were mis-entered for specific accounts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to match the refreshed export timestamp
$ws.Name = "IClientBalance-20240827-094617-"

# Update the reference date for every row of data (row 1 is the header)
$lastRow = 274
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 45531
}

# Fix specific projected value / total value corrections
$ws.Cells.Item(101, 5).Value = 69.98
$ws.Cells.Item(101, 8).Value = 69.98

$ws.Cells.Item(112, 5).Value = 15000.38
$ws.Cells.Item(112, 8).Value = 15000.38

$ws.Cells.Item(120, 5).Value = 14880.43
$ws.Cells.Item(120, 8).Value = 14880.43

$ws.Cells.Item(143, 5).Value = 999.9
$ws.Cells.Item(143, 8).Value = 999.9
